# revisi P3P Agustus 2020
# Add the second half of the Agustus 2020 reading table (dates 11-20, Selasa
# onward) starting at F4:I6 on Sheet1, and update the saved view state
# (zoom + selection) to match where the editor left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 4: Tgl=11, Hari=Selasa, Pagi=Ams 11, Sore=Mzm 74
$ws.Range("F4").Value = 11
$ws.Range("G4").Value = "Selasa"
$ws.Range("H4").Value = "Ams 11"
$ws.Range("I4").Value = "Mzm 74"

# Row 5: Pagi=Luk 2, Sore=Hak 12
$ws.Range("H5").Value = "Luk 2"
$ws.Range("I5").Value = "Hak 12"

# Row 6: Pagi=1 Kor 9, Sore=Yeh 5
$ws.Range("H6").Value = "1 Kor 9"
$ws.Range("I6").Value = "Yeh 5"

# Update view: zoom to 115% and move the selection to H7
$ws.Range("H7").Select() | Out-Null
$ws.Application.ActiveWindow.Zoom = 115
